# "set descriptions for each tests which will be appeared in extent report"
#
# Runner sheet: drop the loginTest/logoutTest/loginTestInValid rows entirely,
# and give each remaining UI-header test (bestSellersTest / newReleaseTest /
# moversandShakersTest) its own specific description instead of the one
# generic sentence they all used to share. CredentialData's own data is
# untouched, it just stops being the active tab (Runner becomes active).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Runner")
$ws2 = $wb.Worksheets.Item("CredentialData")

# --- Runner: move the bestSellersTest/newReleaseTest/moversandShakersTest
# rows (5-7) up onto rows 2-4, overwriting the login/logout/invalid-login
# rows while preserving cell styling, then drop the now-duplicate tail ---
$ws1.Range("A5:E7").Copy($ws1.Range("A2"))
$ws1.Range("A5:A7").EntireRow.Delete() | Out-Null

# --- Runner: give each remaining test its own specific description, and
# flip bestSellersTest to not-executed. Values are staged through a scratch
# cell and pasted as values-only so the destination cell keeps its original
# style (incl. quotePrefix) instead of forking a new one. ---
$scratch = $ws1.Range("ZZ100")

function Set-CellText($range, $text) {
    $scratch.Value = $text
    $scratch.Copy()
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

Set-CellText $ws1.Range("B2") "To verity header of the best sellers page is correct or not"
Set-CellText $ws1.Range("C2") "No"

Set-CellText $ws1.Range("B3") "To verity header of the new release page is correct or not"
Set-CellText $ws1.Range("C3") "Yes"

Set-CellText $ws1.Range("B4") "To verity header of the mover and shaker page is correct or not"
Set-CellText $ws1.Range("C4") "Yes"

$scratch.ClearContents() | Out-Null
$excel.CutCopyMode = $false

# --- Runner: widen the description column to fit the new text ---
$ws1.Columns.Item(2).ColumnWidth = 52.5

# --- Runner becomes the active sheet/tab with B14 selected; CredentialData
# keeps its own F9 selection but is no longer the active tab ---
$ws2.Range("F9").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("B14").Select() | Out-Null
